$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (prices/volumes are text-formatted cells,
# so we prefix with a leading apostrophe to force text entry and avoid
# Excel auto-converting them to numeric values, which would lose exact
# formatting such as trailing zeros, multiple dot separators, and padding.

$ws.Range("D2").Value = "'30.460.76"
$ws.Range("E2").Value = "'  +0.53%  "

$ws.Range("D3").Value = "'2.106.95"
$ws.Range("E3").Value = "'  +1.13%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "'  +0.65%  "

$ws.Range("E5").Value = "'  +2.03%  "

$ws.Range("E6").Value = "'  +0.71%  "

$ws.Range("D7").Value = "'0.5228"
$ws.Range("E7").Value = "'  +0.65%  "

$ws.Range("D8").Value = "'0.4570"
$ws.Range("E8").Value = "'  +5.85%  "

$ws.Range("D9").Value = "'53.31"
$ws.Range("E9").Value = "'  +15.60%  "

$ws.Range("D10").Value = "'0.08952"
$ws.Range("E10").Value = "'  +1.37%  "

$ws.Range("D11").Value = "'1.176"
$ws.Range("E11").Value = "'  +1.44%  "

$ws.Range("D12").Value = "'24.26"
$ws.Range("E12").Value = "'  -0.64%  "

$ws.Range("D13").Value = "'2.100.69"
$ws.Range("E13").Value = "'  +1.52%  "

$ws.Range("D14").Value = "'6.855"
$ws.Range("E14").Value = "'  +2.75%  "

$ws.Range("D15").Value = "'8.041"
$ws.Range("E15").Value = "'  +4.63%  "

$ws.Range("D16").Value = "'96.65"
$ws.Range("E16").Value = "'  +1.62%  "

$ws.Range("D17").Value = "'0.00001147"
$ws.Range("E17").Value = "'  +2.39%  "

$ws.Range("E18").Value = "'  +0.76%  "

$ws.Range("D19").Value = "'0.06659"
$ws.Range("E19").Value = "'  +0.74%  "

$ws.Range("D20").Value = "'19.23"
$ws.Range("E20").Value = "'  +2.23%  "

$ws.Range("E21").Value = "'  +0.73%  "

$ws.Range("D22").Value = "'6.343"
$ws.Range("E22").Value = "'  +0.53%  "

$ws.Range("D23").Value = "'30.512.16"
$ws.Range("E23").Value = "'  +0.71%  "

$ws.Range("D24").Value = "'12.51"
$ws.Range("E24").Value = "'  +1.95%  "

$ws.Range("D25").Value = "'2.360"
$ws.Range("E25").Value = "'  +3.17%  "

$ws.Range("D26").Value = "'2.345.43"
$ws.Range("E26").Value = "'  +1.40%  "

$ws.Range("D27").Value = "'22.30"
$ws.Range("E27").Value = "'  -0.35%  "

$ws.Range("B28").Value = "'LidoDAOToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.543"
$ws.Range("E28").Value = "'  -1.51%  "

$ws.Range("B29").Value = "'Monero"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'163.22"
$ws.Range("E29").Value = "'  +1.00%  "

$ws.Range("D30").Value = "'133.15"
$ws.Range("E30").Value = "'  +1.71%  "

$ws.Range("D31").Value = "'1.218"
$ws.Range("E31").Value = "'  +2.64%  "

$ws.Range("D32").Value = "'0.1073"
$ws.Range("E32").Value = "'  +0.75%  "

$ws.Range("D33").Value = "'1.661"
$ws.Range("E33").Value = "'  +1.20%  "

$ws.Range("D34").Value = "'6.377"
$ws.Range("E34").Value = "'  +3.00%  "

$ws.Range("E35").Value = "'  +3.31%  "

$ws.Range("D36").Value = "'10.37"
$ws.Range("E36").Value = "'  +5.77%  "

$ws.Range("B37").Value = "'InternetComputer(DFINITY)"
$ws.Range("C37").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.722"
$ws.Range("E37").Value = "'  +5.37%  "

$ws.Range("B38").Value = "'VeChain"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02580"
$ws.Range("E38").Value = "'  +0.11%  "

$ws.Range("D39").Value = "'0.06842"
$ws.Range("E39").Value = "'  +3.06%  "

$ws.Range("E40").Value = "'  +2.58%  "

$ws.Range("D41").Value = "'12.71"
$ws.Range("E41").Value = "'  +0.20%  "

$ws.Range("D42").Value = "'0.6881"
$ws.Range("E42").Value = "'  +1.21%  "

$ws.Range("D43").Value = "'1.254"
$ws.Range("E43").Value = "'  +0.36%  "

$ws.Range("E44").Value = "'  +5.75%  "

$ws.Range("E45").Value = "'  +0.54%  "

$ws.Range("D46").Value = "'0.6389"
$ws.Range("E46").Value = "'  +0.75%  "

$ws.Range("D47").Value = "'3.666"
$ws.Range("E47").Value = "'  +1.74%  "

$ws.Range("D48").Value = "'1.254"
$ws.Range("E48").Value = "'  +1.36%  "

$ws.Range("D49").Value = "'0.3427"
$ws.Range("E49").Value = "'  +25.51%  "

$ws.Range("B50").Value = "'Aave"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'83.31"
$ws.Range("E50").Value = "'  +2.22%  "

$ws.Range("D51").Value = "'1.205"
$ws.Range("E51").Value = "'  +1.09%  "
